$d = $word.ActiveDocument

# The three header images ("image1.png","image2.png","image3.png") are
# identified in the header's drawing XML (wp:docPr / pic:cNvPr name
# attributes). The edit swaps the display names used for the 923925x695325
# logo (docPr id="3") and the 5467350x38100 divider line (docPr id="2"):
#   docPr id="3" name="image3.png" -> name="image2.png"
#   docPr id="2" name="image2.png" -> name="image3.png"
# (and the matching pic:cNvPr name= for each), while the embedded
# relationship ids (and therefore the actual picture bytes) stay the same.
#
# The document has no section properties wiring header1.xml into the
# Sections/Headers object model, so Find/Replace on the header Range can't
# reach it. Document.WordOpenXML round-trips the full package (including
# otherwise-unreachable parts like this header), so we edit the raw OOXML
# there instead.

$xml = $d.WordOpenXML

$markerA = "__IMG_SWAP_MARKER_A__"
$markerB = "__IMG_SWAP_MARKER_B__"

# docPr id="3" ... name="image3.png"  <->  name="image2.png"
$xml = $xml.Replace('<wp:docPr id="3" name="image3.png"/>', '<wp:docPr id="3" name="' + $markerA + '"/>')
$xml = $xml.Replace('<wp:docPr id="2" name="image2.png"/>', '<wp:docPr id="2" name="image3.png"/>')
$xml = $xml.Replace('<wp:docPr id="3" name="' + $markerA + '"/>', '<wp:docPr id="3" name="image2.png"/>')

# pic:cNvPr id="0" name="image3.png" (inside the id="3" picture) <-> name="image2.png"
$xml = $xml.Replace('<pic:cNvPr id="0" name="image3.png"/><pic:cNvPicPr preferRelativeResize="0"/></pic:nvPicPr><pic:blipFill><a:blip r:embed="rId2"/>', '<pic:cNvPr id="0" name="' + $markerB + '"/><pic:cNvPicPr preferRelativeResize="0"/></pic:nvPicPr><pic:blipFill><a:blip r:embed="rId2"/>')
$xml = $xml.Replace('<pic:cNvPr id="0" name="image2.png"/><pic:cNvPicPr preferRelativeResize="0"/></pic:nvPicPr><pic:blipFill><a:blip r:embed="rId3"/>', '<pic:cNvPr id="0" name="image3.png"/><pic:cNvPicPr preferRelativeResize="0"/></pic:nvPicPr><pic:blipFill><a:blip r:embed="rId3"/>')
$xml = $xml.Replace('<pic:cNvPr id="0" name="' + $markerB + '"/><pic:cNvPicPr preferRelativeResize="0"/></pic:nvPicPr><pic:blipFill><a:blip r:embed="rId2"/>', '<pic:cNvPr id="0" name="image2.png"/><pic:cNvPicPr preferRelativeResize="0"/></pic:nvPicPr><pic:blipFill><a:blip r:embed="rId2"/>')

$d.WordOpenXML = $xml

Write-Host "Header image names swapped"
